$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.351589679718018
$ws.Range("B1").Value = 1.532183408737183
$ws.Range("C1").Value = 4.005435466766357
$ws.Range("D1").Value = 3.193935871124268
$ws.Range("E1").Value = 1.083919286727905
